# Generate Report for Handback
#
# This mirrors a localization "handback" run: the zh-cn and de-de targets
# have come back in sync with en-US, so the Overview / per-locale status
# sheets get new status text, the handback file name + datetime are filled
# in, and a couple of columns get widened so the longer strings/dates are
# readable.

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"
$targetFileName = "c26435de-ab30-4a72-80d7-c8b73e003e4d.md"
$targetDisplay  = "c26435de-ab30-4a72-80d7-c8b73e003e4d.md"

# ---------------------------------------------------------------------------
# Overview sheet: zh-cn / de-de status columns (E2, F2)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusNew
$wsOverview.Range("F2").Value = $statusNew

# Widen the two status columns now that the text is longer.
$wsOverview.Range("E1").EntireColumn.ColumnWidth = 29.17
$wsOverview.Range("F1").EntireColumn.ColumnWidth = 29.17

# ---------------------------------------------------------------------------
# zh-cn status sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Status
$wsZh.Range("C2").Value = $statusNew

# Latest Target File (I2) - becomes a hyperlink to the source .md file,
# same as A2 / "Source File Name".
$wsZh.Range("I2").Value = $targetFileName
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bd29aa27ade1b82f76aecc5272d44d29dcc51171/e2e/c26435de-ab30-4a72-80d7-c8b73e003e4d.md", "", "", $targetDisplay) | Out-Null
$wsZh.Range("I2").Style = "HyperLink"

# Latest Handback File (J2)
$wsZh.Range("J2").Value = "c26435de-ab30-4a72-80d7-c8b73e003e4d.4704f6817069b3df45e2962c1a742c0cd81e8083.zh-cn.xlf"

# Latest Handback DateTime (K2)
$wsZh.Range("K2").Value = "2016-08-25 04:57:07"

# Column widths: Status (C), Latest Target File (I), Latest Handback File (J)
$wsZh.Range("C1").EntireColumn.ColumnWidth = 29.17
$wsZh.Range("I1").EntireColumn.ColumnWidth = 39.17
$wsZh.Range("J1").EntireColumn.ColumnWidth = 39.17

# ---------------------------------------------------------------------------
# de-de status sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Status
$wsDe.Range("C2").Value = $statusNew

# Latest Target File (I2) - becomes a hyperlink to the source .md file,
# same as A2 / "Source File Name".
$wsDe.Range("I2").Value = $targetFileName
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bd29aa27ade1b82f76aecc5272d44d29dcc51171/e2e/c26435de-ab30-4a72-80d7-c8b73e003e4d.md", "", "", $targetDisplay) | Out-Null
$wsDe.Range("I2").Style = "HyperLink"

# Latest Handback File (J2)
$wsDe.Range("J2").Value = "c26435de-ab30-4a72-80d7-c8b73e003e4d.4704f6817069b3df45e2962c1a742c0cd81e8083.de-de.xlf"

# Latest Handback DateTime (K2)
$wsDe.Range("K2").Value = "2016-08-25 04:57:14"

# Column widths: Status (C), Latest Target File (I), Latest Handback File (J)
$wsDe.Range("C1").EntireColumn.ColumnWidth = 29.17
$wsDe.Range("I1").EntireColumn.ColumnWidth = 39.17
$wsDe.Range("J1").EntireColumn.ColumnWidth = 39.17
